$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3023.7778
$ws.Range("I137").Value = 1399.8572
$ws.Range("J137").Value = 4057.182
$ws.Range("K137").Value = 4199.571599999999
$ws.Range("L137").Value = 12171.546
$ws.Range("M137").Value = -1649.571599999999
$ws.Range("N137").Value = -17271.546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 1127
$ws.Range("I10").Value = 1127
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 1127
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -957
$ws.Range("N10").Value = ""

$ws.Range("H32").Value = 7008.421
$ws.Range("I32").Value = 7008.421
$ws.Range("K32").Value = 7008.421
$ws.Range("M32").Value = -6721.421

$ws.Range("H63").Value = 2938.4
$ws.Range("I63").Value = 2230.6667
$ws.Range("K63").Value = 2230.6667
$ws.Range("M63").Value = -1544.6667

$ws.Range("H66").Value = 2938.4
$ws.Range("I66").Value = 2230.6667
$ws.Range("K66").Value = 11153.3335
$ws.Range("M66").Value = -7721.333500000001

$ws.Range("H88").Value = 4234.5
$ws.Range("I88").Value = 1540
$ws.Range("J88").Value = 5132.6665
$ws.Range("K88").Value = 1540
$ws.Range("L88").Value = 5132.6665
$ws.Range("M88").Value = -1134
$ws.Range("N88").Value = -5944.6665

$ws.Range("H91").Value = 4234.5
$ws.Range("I91").Value = 1540
$ws.Range("J91").Value = 5132.6665
$ws.Range("K91").Value = 1540
$ws.Range("L91").Value = 5132.6665
$ws.Range("M91").Value = -136
$ws.Range("N91").Value = -7940.6665

$ws.Range("H114").Value = 20000
$ws.Range("J114").Value = 20000
$ws.Range("L114").Value = 20000
$ws.Range("N114").Value = -28678

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = ""

$ws.Range("H20").Value = 4473.8423
$ws.Range("I20").Value = 4008.4285
$ws.Range("J20").Value = 5777
$ws.Range("K20").Value = 4008.4285
$ws.Range("L20").Value = 5777
$ws.Range("M20").Value = -3761.4285
$ws.Range("N20").Value = -6271

$ws.Range("H62").Value = 35000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = ""

$ws.Range("H65").Value = 35000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = ""

$ws.Range("H94").Value = 2169.5
$ws.Range("I94").Value = 2281.739
$ws.Range("J94").Value = 1309
$ws.Range("K94").Value = 2281.739
$ws.Range("L94").Value = 1309
$ws.Range("M94").Value = -1830.739
$ws.Range("N94").Value = -2211

$ws.Range("H99").Value = 3400
$ws.Range("I99").Value = 3400
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3400
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1902
$ws.Range("N99").Value = ""

$ws.Range("H134").Value = 3347.6538
$ws.Range("I134").Value = 3385.6
$ws.Range("K134").Value = 10156.8
$ws.Range("M134").Value = -7621.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 775
$ws.Range("I22").Value = 750
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 750
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = -400
$ws.Range("N22").Value = -1500

$ws.Range("H31").Value = 2424.7727
$ws.Range("I31").Value = 1976.75
$ws.Range("K31").Value = 1976.75
$ws.Range("M31").Value = -1681.75

$ws.Range("H34").Value = 2424.7727
$ws.Range("I34").Value = 1976.75
$ws.Range("K34").Value = 1976.75
$ws.Range("M34").Value = -1774.75

$ws.Range("H86").Value = 13777.667
$ws.Range("I86").Value = 13166.5
$ws.Range("K86").Value = 13166.5
$ws.Range("M86").Value = -12043.5

$ws.Range("H89").Value = 13777.667
$ws.Range("I89").Value = 13166.5
$ws.Range("K89").Value = 65832.5
$ws.Range("M89").Value = -60216.5

$ws.Range("H122").Value = 625.1111
$ws.Range("I122").Value = 625.1111
$ws.Range("K122").Value = 1875.3333
$ws.Range("M122").Value = 574.6667000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2251.8572
$ws.Range("I75").Value = 750
$ws.Range("J75").Value = 2502.1667
$ws.Range("K75").Value = 2250
$ws.Range("L75").Value = 7506.500100000001
$ws.Range("M75").Value = -1252
$ws.Range("N75").Value = -9502.500100000001

$ws.Range("H78").Value = 2251.8572
$ws.Range("I78").Value = 750
$ws.Range("J78").Value = 2502.1667
$ws.Range("K78").Value = 6750
$ws.Range("L78").Value = 22519.5003
$ws.Range("M78").Value = -1758
$ws.Range("N78").Value = -32503.5003

$ws.Range("H92").Value = 493
$ws.Range("I92").Value = 416.25
$ws.Range("J92").Value = 800
$ws.Range("K92").Value = 1248.75
$ws.Range("L92").Value = 2400
$ws.Range("M92").Value = -0.75
$ws.Range("N92").Value = -4896

$ws.Range("H114").Value = 199.5
$ws.Range("I114").Value = 199
$ws.Range("J114").Value = 200
$ws.Range("K114").Value = 597
$ws.Range("L114").Value = 600
$ws.Range("M114").Value = 2657
$ws.Range("N114").Value = -7108

$ws.Range("H120").Value = 4353.6665
$ws.Range("I120").Value = 363
$ws.Range("K120").Value = 1089
$ws.Range("M120").Value = 3749

$ws.Range("H122").Value = 477.46667
$ws.Range("I122").Value = 384.66666
$ws.Range("J122").Value = 500.66666
$ws.Range("K122").Value = 3461.99994
$ws.Range("L122").Value = 4505.99994
$ws.Range("M122").Value = -1011.99994
$ws.Range("N122").Value = -9405.99994

$ws.Range("H129").Value = 543.4286
$ws.Range("I129").Value = 509
$ws.Range("K129").Value = 1527
$ws.Range("M129").Value = 3473

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 232.11111
$ws.Range("I2").Value = 297.5
$ws.Range("K2").Value = 297.5
$ws.Range("M2").Value = -184.5

$ws.Range("H12").Value = 3003
$ws.Range("I12").Value = 3003
$ws.Range("K12").Value = 3003
$ws.Range("M12").Value = -2863

$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").Value = ""

$ws.Range("H57").Value = 22998.5
$ws.Range("I57").Value = 12000
$ws.Range("J57").Value = 25198.2
$ws.Range("K57").Value = 12000
$ws.Range("L57").Value = 25198.2
$ws.Range("M57").Value = -11180
$ws.Range("N57").Value = -26838.2

$ws.Range("H94").Value = 50000
$ws.Range("J94").Value = 50000
$ws.Range("L94").Value = 50000
$ws.Range("N94").Value = -51352

$ws.Range("H97").Value = 772.44446
$ws.Range("I97").Value = 550.2857
$ws.Range("J97").Value = 1550
$ws.Range("K97").Value = 550.2857
$ws.Range("L97").Value = 1550
$ws.Range("M97").Value = -54.28570000000002
$ws.Range("N97").Value = -2542

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1081.1818
$ws.Range("I22").Value = 1127.5714
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 1127.5714
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -832.5714
$ws.Range("N22").Value = -1590

$ws.Range("H27").Value = 1081.1818
$ws.Range("I27").Value = 1127.5714
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 1127.5714
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -1020.5714
$ws.Range("N27").Value = -1214

$ws.Range("H46").Value = 3062.6453
$ws.Range("J46").Value = 3765.1875
$ws.Range("L46").Value = 3765.1875
$ws.Range("N46").Value = -4141.1875

$ws.Range("H87").Value = 39997
$ws.Range("J87").Value = 39997
$ws.Range("L87").Value = 39997
$ws.Range("N87").Value = -42243

$ws.Range("H90").Value = 39997
$ws.Range("J90").Value = 39997
$ws.Range("L90").Value = 119991
$ws.Range("N90").Value = -131223

$ws.Range("H93").Value = 839
$ws.Range("I93").Value = 850
$ws.Range("J93").Value = 795
$ws.Range("K93").Value = 850
$ws.Range("L93").Value = 795
$ws.Range("M93").Value = 398
$ws.Range("N93").Value = -3291

$ws.Range("H132").Value = 4706
$ws.Range("I132").Value = 3590.5
$ws.Range("K132").Value = 10771.5
$ws.Range("M132").Value = -8241.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8899.799999999999
$ws.Range("I136").Value = 8524.75
$ws.Range("K136").Value = 25574.25
$ws.Range("M136").Value = -23024.25
